$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column layout: A Sending cluster, B Ligand symbol, C Receptor symbol, D Target cluster,
# E..T numeric metrics.

# ---- Row 2 ----
$ws.Cells.Item(2,1).Value  = "FAPs"
$ws.Cells.Item(2,2).Value  = "Ntn3"
$ws.Cells.Item(2,3).Value  = "Cdon"
$ws.Cells.Item(2,4).Value  = "ECs"
$ws.Cells.Item(2,5).Value  = 3
$ws.Cells.Item(2,6).Value  = 1
$ws.Cells.Item(2,7).Value  = 2.147752666666667
$ws.Cells.Item(2,8).Value  = 6.443258
$ws.Cells.Item(2,9).Value  = 0.8708065646157738
$ws.Cells.Item(2,10).Value = 0.8708065646157739
$ws.Cells.Item(2,11).Value = 2
$ws.Cells.Item(2,12).Value = 0.6666666666666666
$ws.Cells.Item(2,13).Value = 0.9894223333333333
$ws.Cells.Item(2,14).Value = 2.968267
$ws.Cells.Item(2,15).Value = 0.02008622060203138
$ws.Cells.Item(2,16).Value = 0.02008622060203139
$ws.Cells.Item(2,17).Value = 2.125034454876222
$ws.Cells.Item(2,18).Value = 19.125310093886
$ws.Cells.Item(2,19).Value = 0.01749121275856953
$ws.Cells.Item(2,20).Value = 0.01749121275856954

# ---- Row 3 ----
$ws.Cells.Item(3,1).Value  = "FAPs"
$ws.Cells.Item(3,2).Value  = "Ntn3"
$ws.Cells.Item(3,3).Value  = "Cdon"
$ws.Cells.Item(3,4).Value  = "FAPs"
$ws.Cells.Item(3,5).Value  = 3
$ws.Cells.Item(3,6).Value  = 1
$ws.Cells.Item(3,7).Value  = 2.147752666666667
$ws.Cells.Item(3,8).Value  = 6.443258
$ws.Cells.Item(3,9).Value  = 0.8708065646157738
$ws.Cells.Item(3,10).Value = 0.8708065646157739
$ws.Cells.Item(3,11).Value = 3
$ws.Cells.Item(3,12).Value = 1
$ws.Cells.Item(3,13).Value = 34.79912266666667
$ws.Cells.Item(3,14).Value = 104.397368
$ws.Cells.Item(3,15).Value = 0.7064555054917406
$ws.Cells.Item(3,16).Value = 0.7064555054917406
$ws.Cells.Item(3,17).Value = 74.73990850499379
$ws.Cells.Item(3,18).Value = 672.659176544944
$ws.Cells.Item(3,19).Value = 0.6151860917911626
$ws.Cells.Item(3,20).Value = 0.6151860917911627

# ---- Row 4 ----
$ws.Cells.Item(4,1).Value  = "FAPs"
$ws.Cells.Item(4,2).Value  = "Ntn3"
$ws.Cells.Item(4,3).Value  = "Cdon"
$ws.Cells.Item(4,4).Value  = "sCs"
$ws.Cells.Item(4,5).Value  = 3
$ws.Cells.Item(4,6).Value  = 1
$ws.Cells.Item(4,7).Value  = 2.147752666666667
$ws.Cells.Item(4,8).Value  = 6.443258
$ws.Cells.Item(4,9).Value  = 0.8708065646157738
$ws.Cells.Item(4,10).Value = 0.8708065646157739
$ws.Cells.Item(4,11).Value = 3
$ws.Cells.Item(4,12).Value = 1
$ws.Cells.Item(4,13).Value = 13.47021566666667
$ws.Cells.Item(4,14).Value = 40.410647
$ws.Cells.Item(4,15).Value = 0.2734582739062281
$ws.Cells.Item(4,16).Value = 0.2734582739062281
$ws.Cells.Item(4,17).Value = 28.93069161865845
$ws.Cells.Item(4,18).Value = 260.376224567926
$ws.Cells.Item(4,19).Value = 0.2381292600660418
$ws.Cells.Item(4,20).Value = 0.2381292600660418

# ---- Row 5 ----
$ws.Cells.Item(5,1).Value  = "sCs"
$ws.Cells.Item(5,2).Value  = "Ntn3"
$ws.Cells.Item(5,3).Value  = "Cdon"
$ws.Cells.Item(5,4).Value  = "ECs"
$ws.Cells.Item(5,5).Value  = 3
$ws.Cells.Item(5,6).Value  = 1
$ws.Cells.Item(5,7).Value  = 0.318642
$ws.Cells.Item(5,8).Value  = 0.9559260000000001
$ws.Cells.Item(5,9).Value  = 0.1291934353842261
$ws.Cells.Item(5,10).Value = 0.1291934353842261
$ws.Cells.Item(5,11).Value = 2
$ws.Cells.Item(5,12).Value = 0.6666666666666666
$ws.Cells.Item(5,13).Value = 0.9894223333333333
$ws.Cells.Item(5,14).Value = 2.968267
$ws.Cells.Item(5,15).Value = 0.02008622060203138
$ws.Cells.Item(5,16).Value = 0.02008622060203139
$ws.Cells.Item(5,17).Value = 0.315271511138
$ws.Cells.Item(5,18).Value = 2.837443600242
$ws.Cells.Item(5,19).Value = 0.002595007843461854
$ws.Cells.Item(5,20).Value = 0.002595007843461854

# ---- Row 6 (new) ----
$ws.Cells.Item(6,1).Value  = "sCs"
$ws.Cells.Item(6,2).Value  = "Ntn3"
$ws.Cells.Item(6,3).Value  = "Cdon"
$ws.Cells.Item(6,4).Value  = "FAPs"
$ws.Cells.Item(6,5).Value  = 3
$ws.Cells.Item(6,6).Value  = 1
$ws.Cells.Item(6,7).Value  = 0.318642
$ws.Cells.Item(6,8).Value  = 0.9559260000000001
$ws.Cells.Item(6,9).Value  = 0.1291934353842261
$ws.Cells.Item(6,10).Value = 0.1291934353842261
$ws.Cells.Item(6,11).Value = 3
$ws.Cells.Item(6,12).Value = 1
$ws.Cells.Item(6,13).Value = 34.79912266666667
$ws.Cells.Item(6,14).Value = 104.397368
$ws.Cells.Item(6,15).Value = 0.7064555054917406
$ws.Cells.Item(6,16).Value = 0.7064555054917406
$ws.Cells.Item(6,17).Value = 11.088462044752
$ws.Cells.Item(6,18).Value = 99.796158402768
$ws.Cells.Item(6,19).Value = 0.09126941370057801
$ws.Cells.Item(6,20).Value = 0.09126941370057801

# ---- Row 7 (new) ----
$ws.Cells.Item(7,1).Value  = "sCs"
$ws.Cells.Item(7,2).Value  = "Ntn3"
$ws.Cells.Item(7,3).Value  = "Cdon"
$ws.Cells.Item(7,4).Value  = "sCs"
$ws.Cells.Item(7,5).Value  = 3
$ws.Cells.Item(7,6).Value  = 1
$ws.Cells.Item(7,7).Value  = 0.318642
$ws.Cells.Item(7,8).Value  = 0.9559260000000001
$ws.Cells.Item(7,9).Value  = 0.1291934353842261
$ws.Cells.Item(7,10).Value = 0.1291934353842261
$ws.Cells.Item(7,11).Value = 3
$ws.Cells.Item(7,12).Value = 1
$ws.Cells.Item(7,13).Value = 13.47021566666667
$ws.Cells.Item(7,14).Value = 40.410647
$ws.Cells.Item(7,15).Value = 0.2734582739062281
$ws.Cells.Item(7,16).Value = 0.2734582739062281
$ws.Cells.Item(7,17).Value = 4.292176460458
$ws.Cells.Item(7,18).Value = 38.629588144122
$ws.Cells.Item(7,19).Value = 0.03532901384018629
$ws.Cells.Item(7,20).Value = 0.03532901384018629
